$wb = $excel.ActiveWorkbook

# Revert "Temporarily disabling feedback loops until issues are sorted out."
# BDMFL!B2 goes back to 0 (feedback loops enabled / not disabled).
$ws = $wb.Worksheets.Item("BDMFL")
$ws.Range("B2").Value = 0

# Restore the active sheet/selection to the "About" sheet (tab 1),
# matching the pre-revert workbook view state.
$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Activate()
$aboutWs.Range("A1").Select()
